$wb = $excel.ActiveWorkbook

# The two language sheets
$wsEn = $wb.Worksheets.Item(1)   # "en"
$wsEs = $wb.Worksheets.Item(2)   # "es"

# New key/value pairs for the act 3 intro dialogue
$keys = @("act_3_intro_1", "act_3_intro_2", "act_3_intro_3", "act_3_intro_4")

$rightQuote = [char]0x2019
$values = @(
    "In this act, we will be guiding the golden ball to its holy source.",
    "To do that, you must place force fields to influence the golden ball" + $rightQuote + "s movement.",
    "Remember that so long as any forces are acting on an object, its acceleration will change. Thus, causing the object" + $rightQuote + "s speed to change over time. ",
    "Now go forth, and attach a force field on the wall to get the ball rolling!"
)

$startRow = 123

# First write all the keys (column A) for both sheets, row by row.
for ($i = 0; $i -lt $keys.Length; $i++) {
    $row = $startRow + $i
    $wsEn.Cells.Item($row, 1).Value = $keys[$i]
}

# Then write the English values (column B), row by row.
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $wsEn.Cells.Item($row, 2).Value = $values[$i]
}

# Spanish sheet: only the key in column A (translation not filled in yet)
for ($i = 0; $i -lt $keys.Length; $i++) {
    $row = $startRow + $i
    $wsEs.Cells.Item($row, 1).Value = $keys[$i]
}

# Update selection to match where the editor ended up on each sheet
$wsEn.Activate()
$wsEn.Range("B122").Select()

$wsEs.Activate()
$wsEs.Range("B123").Select()

$wsEn.Activate()
